# Update the imputed values in the RandomForest result sheet.
# (Commit message: "Update Name of Algo" - underlying repo path/algo name
#  changed; the data values below reflect the refreshed run for this algo.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.3963
$ws.Range("A12").Value = -21.605
$ws.Range("C14").Value = -13.1976
$ws.Range("C26").Value = -12.52320000000001
$ws.Range("C31").Value = -13.0687
$ws.Range("A32").Value = -21.22859999999999
$ws.Range("C35").Value = -12.43720000000001
$ws.Range("A36").Value = -20.1582
$ws.Range("C37").Value = -13.5607
$ws.Range("A38").Value = -19.42449999999999
$ws.Range("C45").Value = -13.85759999999999
$ws.Range("A46").Value = -21.71669999999999
$ws.Range("A54").Value = -21.76439999999999
$ws.Range("A55").Value = -22.44110000000001
$ws.Range("C57").Value = -14.31559999999999
$ws.Range("A67").Value = -21.54929999999997
$ws.Range("A69").Value = -21.64349999999997
$ws.Range("A72").Value = -21.71919999999998
$ws.Range("A91").Value = -21.42750000000001
$ws.Range("A99").Value = -20.51299999999999
$ws.Range("C100").Value = -13.11809999999999
$ws.Range("C102").Value = -12.9886
